$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new "2022" column (S) ------------------------------------
# Copy formatting from the preceding column (R) so the new cells inherit
# the same visual style, then set the actual values.
$ws.Range("R4").Copy($ws.Range("S4"))
$ws.Range("S4").Value = 2022

$ws.Range("R5").Copy($ws.Range("S5"))
$ws.Range("S5").Value = 3.4

# --- Update the existing trend values in row 5 -------------------------
$ws.Range("P5").Value = 4.4
$ws.Range("Q5").Value = 2.9
$ws.Range("R5").Value = 3.2

# --- Update the selected cell shown in the saved view -------------------
$ws.Range("T4").Select()
